$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set column H width to match the other data columns (15)
$ws.Columns.Item(8).ColumnWidth = 14.1666666666667

$ws.Range("B1").Value = "Data collected @ Sat Jul 15 10:43:16 2023.txt"
$ws.Range("C1").Value = "Data collected @ Sat Jul 15 10:43:16 2023.txt"
$ws.Range("D1").Value = "Data collected @ Sat Jul 15 11:06:46 2023.txt"
$ws.Range("E1").Value = "Data collected @ Sat Aug 3 8:10:15 2023.txt"
$ws.Range("F1").Value = "Data collected @ Sat Jul 89 10:55:15 2023.txt"
$ws.Range("G1").Value = "Data collected @ Sat Jul 15 10:45:05 2023.txt"
$ws.Range("H1").Value = "Data collected @ Sat Jul 15 10:55:15 2023.txt"

$ws.Range("B2").Value = "10:43:21"
$ws.Range("C2").Value = "10:43:21"
$ws.Range("D2").Value = "11:07:01"
$ws.Range("E2").Value = "12:55:28"
$ws.Range("F2").Value = "8:25:22"
$ws.Range("G2").Value = "10:45:10"
$ws.Range("H2").Value = "10:55:22"

$ws.Range("B3").Value = "10:43:51"
$ws.Range("C3").Value = "10:43:51"
$ws.Range("D3").Value = "11:11:39"
$ws.Range("E3").Value = "22:55:28"
$ws.Range("F3").Value = "22:55:28"
$ws.Range("G3").Value = "10:50:49"
$ws.Range("H3").Value = "10:57:15"

$ws.Range("B4").Value = "0:0:30"
$ws.Range("C4").Value = "0:0:30"
$ws.Range("D4").Value = "0:4:38"
$ws.Range("E4").Value = "10:0:0"
$ws.Range("F4").Value = "14:30:6"
$ws.Range("G4").Value = "0:5:39"
$ws.Range("H4").Value = "0:1:53"

$ws.Range("B5").Value = 19
$ws.Range("C5").Value = 19
$ws.Range("D5").Value = 70
$ws.Range("E5").Value = 9
$ws.Range("F5").Value = 12
$ws.Range("G5").Value = 103
$ws.Range("H5").Value = 42

$ws.Range("B6").Value = 11
$ws.Range("C6").Value = 11
$ws.Range("D6").Value = 15
$ws.Range("E6").Value = 0
$ws.Range("F6").Value = 0
$ws.Range("G6").Value = 39
$ws.Range("H6").Value = 20

$ws.Range("B7").Value = 8
$ws.Range("C7").Value = 8
$ws.Range("D7").Value = 55
$ws.Range("E7").Value = 9
$ws.Range("F7").Value = 12
$ws.Range("G7").Value = 64
$ws.Range("H7").Value = 22

$ws.Range("B8").Value = 0
$ws.Range("C8").Value = 0
$ws.Range("D8").Value = 0
$ws.Range("E8").Value = 0
$ws.Range("F8").Value = 0
$ws.Range("G8").Value = 0
$ws.Range("H8").Value = 0

$ws.Range("B9").Value = 0.00833
$ws.Range("C9").Value = 0.00833
$ws.Range("D9").Value = 0.07722
$ws.Range("E9").Value = 10
$ws.Range("F9").Value = 14.50167
$ws.Range("G9").Value = 0.09417
$ws.Range("H9").Value = 0.03139

$ws.Range("B10").Value = 2280.91
$ws.Range("C10").Value = 2280.91
$ws.Range("D10").Value = 906.5
$ws.Range("E10").Value = 0.9
$ws.Range("F10").Value = 0.83
$ws.Range("G10").Value = 1093.77
$ws.Range("H10").Value = 1338.01

$ws.Range("B11").Value = 1320.53
$ws.Range("C11").Value = 1320.53
$ws.Range("D11").Value = 194.25
$ws.Range("E11").Value = 0
$ws.Range("F11").Value = 0
$ws.Range("G11").Value = 414.14
$ws.Range("H11").Value = 637.15

$ws.Range("B12").Value = 960.38
$ws.Range("C12").Value = 960.38
$ws.Range("D12").Value = 712.25
$ws.Range("E12").Value = 0.9
$ws.Range("F12").Value = 0.83
$ws.Range("G12").Value = 679.62
$ws.Range("H12").Value = 700.86

$ws.Range("B13").Value = "N/A"
$ws.Range("C13").Value = "N/A"
$ws.Range("D13").Value = "N/A"
$ws.Range("E13").Value = "N/A"
$ws.Range("F13").Value = "8:25:22"
$ws.Range("G13").Value = "N/A"
$ws.Range("H13").Value = "N/A"

$ws.Range("B14").Value = "N/A"
$ws.Range("C14").Value = "N/A"
$ws.Range("D14").Value = "N/A"
$ws.Range("E14").Value = "N/A"
$ws.Range("F14").Value = "10:55:28"
$ws.Range("G14").Value = "N/A"
$ws.Range("H14").Value = "N/A"

$ws.Range("B15").Value = "N/A"
$ws.Range("C15").Value = "N/A"
$ws.Range("D15").Value = "N/A"
$ws.Range("E15").Value = "N/A"
$ws.Range("F15").Value = 3
$ws.Range("G15").Value = "N/A"
$ws.Range("H15").Value = "N/A"

$ws.Range("B16").Value = "N/A"
$ws.Range("C16").Value = "N/A"
$ws.Range("D16").Value = "N/A"
$ws.Range("E16").Value = "N/A"
$ws.Range("F16").Value = 2.50167
$ws.Range("G16").Value = "N/A"
$ws.Range("H16").Value = "N/A"

$ws.Range("B17").Value = "N/A"
$ws.Range("C17").Value = "N/A"
$ws.Range("D17").Value = "N/A"
$ws.Range("E17").Value = "N/A"
$ws.Range("F17").Value = 1.2
$ws.Range("G17").Value = "N/A"
$ws.Range("H17").Value = "N/A"

$ws.Range("B18").Value = "N/A"
$ws.Range("C18").Value = "N/A"
$ws.Range("D18").Value = "N/A"
$ws.Range("E18").Value = "17:55:28"
$ws.Range("F18").Value = "17:55:28"
$ws.Range("G18").Value = "N/A"
$ws.Range("H18").Value = "N/A"

$ws.Range("B19").Value = "N/A"
$ws.Range("C19").Value = "N/A"
$ws.Range("D19").Value = "N/A"
$ws.Range("E19").Value = "21:12:28"
$ws.Range("F19").Value = "21:12:28"
$ws.Range("G19").Value = "N/A"
$ws.Range("H19").Value = "N/A"

$ws.Range("B20").Value = "N/A"
$ws.Range("C20").Value = "N/A"
$ws.Range("D20").Value = "N/A"
$ws.Range("E20").Value = 6
$ws.Range("F20").Value = 6
$ws.Range("G20").Value = "N/A"
$ws.Range("H20").Value = "N/A"

$ws.Range("B21").Value = "N/A"
$ws.Range("C21").Value = "N/A"
$ws.Range("D21").Value = "N/A"
$ws.Range("E21").Value = 3.28333
$ws.Range("F21").Value = 3.28333
$ws.Range("G21").Value = "N/A"
$ws.Range("H21").Value = "N/A"

$ws.Range("B22").Value = "N/A"
$ws.Range("C22").Value = "N/A"
$ws.Range("D22").Value = "N/A"
$ws.Range("E22").Value = 1.83
$ws.Range("F22").Value = 1.83
$ws.Range("G22").Value = "N/A"
$ws.Range("H22").Value = "N/A"

